# Update column C ("Fitness") values per the commit diff:
#  - Rows 2-26  : 7590 -> 7295
#  - Rows 27-252: 7590 or 7573 -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}

for ($r = 27; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
